$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data table: Sending cluster, Ligand symbol, Receptor symbol, Target cluster, then metrics
$clusters = @("ECs", "FAPs", "sCs")

$data = @{
    "ECs|ECs"   = @(2, 0.6666666666666666, 21.78783233333333, 65.363497, 0.9024488799587679, 0.9024488799587679, 3, 1, 11.08503166666667, 33.255095, 0.4259149152633459, 0.4259149152633459, 241.5188113630238, 2173.669302267215, 0.38436643823714, 0.38436643823714)
    "ECs|FAPs"  = @(2, 0.6666666666666666, 21.78783233333333, 65.363497, 0.9024488799587679, 0.9024488799587679, 2, 0.6666666666666666, 0.312365, 0.937095, 0.01200185227312402, 0.01200185227312402, 6.805756246801666, 61.251806221215, 0.01083105814131136, 0.01083105814131136)
    "ECs|sCs"   = @(2, 0.6666666666666666, 21.78783233333333, 65.363497, 0.9024488799587679, 0.9024488799587679, 3, 1, 14.62900266666666, 43.88700799999999, 0.5620832324635302, 0.5620832324635302, 318.7342573052195, 2868.608315746975, 0.5072513835803166, 0.5072513835803166)
    "FAPs|ECs"  = @(3, 1, 1.359006333333333, 4.077019, 0.05628984676448105, 0.05628984676448104, 3, 1, 11.08503166666667, 33.255095, 0.4259149152633459, 0.4259149152633459, 15.06462824020055, 135.581654161805, 0.02397468531488067, 0.02397468531488067)
    "FAPs|FAPs" = @(3, 1, 1.359006333333333, 4.077019, 0.05628984676448105, 0.05628984676448104, 2, 0.6666666666666666, 0.312365, 0.937095, 0.01200185227312402, 0.01200185227312402, 0.4245060133116667, 3.820554119805, 0.0006755824253440894, 0.0006755824253440893)
    "FAPs|sCs"  = @(3, 1, 1.359006333333333, 4.077019, 0.05628984676448105, 0.05628984676448104, 3, 1, 14.62900266666666, 43.88700799999999, 0.5620832324635302, 0.5620832324635302, 19.88090727435022, 178.928165469152, 0.03163957902425629, 0.03163957902425629)
    "sCs|ECs"   = @(3, 1, 0.9961713333333334, 2.988514, 0.04126127327675106, 0.04126127327675105, 3, 1, 11.08503166666667, 33.255095, 0.4259149152633459, 0.4259149152633459, 11.04259077542556, 99.38331697883001, 0.01757379171132518, 0.01757379171132518)
    "sCs|FAPs"  = @(3, 1, 0.9961713333333334, 2.988514, 0.04126127327675106, 0.04126127327675105, 2, 0.6666666666666666, 0.312365, 0.937095, 0.01200185227312402, 0.01200185227312402, 0.3111690585366667, 2.80052152683, 0.0004952117064685659, 0.0004952117064685658)
    "sCs|sCs"   = @(3, 1, 0.9961713333333334, 2.988514, 0.04126127327675106, 0.04126127327675105, 3, 1, 14.62900266666666, 43.88700799999999, 0.5620832324635302, 0.5620832324635302, 14.57299309179022, 131.156937826112, 0.02319226985895731, 0.02319226985895731)
}

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = "Efna1"
        $ws.Cells.Item($row, 3).Value = "Epha2"
        $ws.Cells.Item($row, 4).Value = $target

        $vals = $data["$sending|$target"]
        for ($i = 0; $i -lt $vals.Length; $i++) {
            $ws.Cells.Item($row, 5 + $i).Value = $vals[$i]
        }
        $row++
    }
}
